$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) values

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.405.38"
$ws.Range("E2").Value = "  +4.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.457.26"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.39"
$ws.Range("E5").Value = "  +2.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.51"
$ws.Range("E6").Value = "  +4.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.519"
$ws.Range("E7").Value = "  +1.13%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.537"
$ws.Range("E9").Value = "  +2.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.15"
$ws.Range("E10").Value = "  +2.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0814"
$ws.Range("E11").Value = "  +1.76%  "
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.41"
$ws.Range("E13").Value = "  -3.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.08"
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.845.77"
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.470.90"
$ws.Range("E16").Value = "  +2.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.845"
$ws.Range("E17").Value = "  +1.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "46.235.47"
$ws.Range("E18").Value = "  +3.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.76"
$ws.Range("E19").Value = "  +2.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.44"
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0936"
$ws.Range("E21").Value = "  +1.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.63"
$ws.Range("E22").Value = "  +2.80%  "
$ws.Range("E23").Value = "  +4.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "248.27"
$ws.Range("E24").Value = "  +2.20%  "
$ws.Range("E25").Value = "  +1.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.21"
$ws.Range("E26").Value = "  +3.73%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.22"
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.77"
$ws.Range("E29").Value = "  +2.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.89"
$ws.Range("E30").Value = "  +5.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "49.53"
$ws.Range("E31").Value = "  +2.20%  "
$ws.Range("E32").Value = "  +3.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.80"
$ws.Range("E33").Value = "  +2.13%  "
$ws.Range("E34").Value = "  +3.26%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0766"
$ws.Range("E37").Value = "  +1.85%  "
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("E39").Value = "  +3.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "122.65"
$ws.Range("E40").Value = "  +1.81%  "
$ws.Range("E41").Value = "  +1.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.24"
$ws.Range("E42").Value = "  +0.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.86"
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0293"
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.977.08"
$ws.Range("E45").Value = "  +1.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.99"
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("E47").Value = "  -3.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.87"
$ws.Range("E48").Value = "  +11.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.14"
$ws.Range("E50").Value = "  +10.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.67"
$ws.Range("E51").Value = "  +4.30%  "
